$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Add the new log row (row 7) for Testmail #5
$ws.Cells.Item(7, 1).Value = "Wil je deze klant bellen?"
$ws.Cells.Item(7, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(7, 3).Value = "Testmail #5: Wil je deze klant bellen?"
$ws.Cells.Item(7, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item(7, 5).Value = "Dank voor je bericht. We pakken dit intern op en houden je op de hoogte."
$ws.Cells.Item(7, 6).Value = "2025-07-31 21:32:10"
$ws.Cells.Item(7, 7).Value = "Ja"
$ws.Cells.Item(7, 8).Value = "Ja"
$ws.Cells.Item(7, 9).Value = "Nee"
$ws.Cells.Item(7, 10).Value = "Nee"

# Extend the existing conditional formatting ranges (D/G/H/I/J) down to row 7,
# preserving the original rules/dxfIds by modifying the applies-to range in place.
$ws.Range("D2:D6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D7"))
$ws.Range("G2:G6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G7"))
$ws.Range("H2:H6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H7"))
$ws.Range("I2:I6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I7"))
$ws.Range("J2:J6").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J7"))

# Update Dashboard summary count for "Intern verzoek / Actie voor medewerker"
$dash.Cells.Item(3, 2).Value = 2
